$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -8
    3  = -2
    4  = 3
    6  = 1
    7  = -2
    8  = 10
    10 = -1
    11 = 2
    12 = -4
    13 = -3
    14 = -1
    15 = -3
    16 = -3
    17 = -2
    18 = -1
    19 = -1
    20 = -4
    21 = 1
    23 = -2
    24 = -1
    25 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
